# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamp values that were refreshed when the
# handback report was regenerated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" for 2ea98654-... .md
# Shared between Overview!G4 and de-de!H4 (Correspond Handoff Datetime)
$wsOverview.Range("G4").Value = "2016-08-17 12:43:21"
$wsDeDe.Range("H4").Value     = "2016-08-17 12:43:21"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H4").Value = "2016-08-17 12:43:16"
$wsZhCn.Range("K4").Value = "2016-08-17 12:43:35"

# de-de sheet: Correspond Handback DateTime
$wsDeDe.Range("K4").Value = "2016-08-17 12:43:43"
